# Local test run finished with a better (lower) buffer count, so update the
# submitted/local score on Sheet1. The "Num Improvement" (F4 = F2-F3) and
# "% Improvement" (F5 = F4/F2) cells are formulas and recalculate on their own.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("F3").Value = 929705

# Reposition the saved workbook window to where the author left it.
$excel.ActiveWindow.Left = 16005
$excel.ActiveWindow.Top = 5100

# Leave the cell cursor where the author left it too.
$ws.Range("G11").Select()
